# "change NLU file to separated with scenarios"
#
# Adds a "scenario" tag (column C) on the Vertexes sheet that classifies
# each entity row as a concept (概念), a behavior (行为) or an attribute
# (属性) -- separating the previously single NLU file into scenarios.

$wb = $excel.ActiveWorkbook

$vertexes = $wb.Worksheets.Item("Vertexes")
$edges    = $wb.Worksheets.Item("Edges")

# New scenario-tag column values for rows 2-14 (row 1 is the header row,
# whose label "实体标签" / C1 is unchanged).
$vertexes.Range("C2").Value  = "概念"
$vertexes.Range("C3").Value  = "行为"
$vertexes.Range("C4").Value  = "行为"
$vertexes.Range("C5").Value  = "属性"
$vertexes.Range("C6").Value  = "概念"
$vertexes.Range("C7").Value  = "概念"
$vertexes.Range("C8").Value  = "概念"
$vertexes.Range("C9").Value  = "行为"
$vertexes.Range("C10").Value = "行为"
$vertexes.Range("C11").Value = "行为"
$vertexes.Range("C12").Value = "属性"
$vertexes.Range("C13").Value = "概念"
$vertexes.Range("C14").Value = "概念"

# Column width tweaks to accommodate the new/longer content.
$vertexes.Columns.Item(1).ColumnWidth = 18.565104166666668
$edges.Columns.Item(2).ColumnWidth = 29.830729166666668

# Restore the active selection to A2 on the Vertexes sheet.
$vertexes.Activate() | Out-Null
$vertexes.Range("A2").Select() | Out-Null
